# Apply the recorded edits to the presentation.
#
# 1) Slide 5 contains a 3x6 table whose table style is switched from the
#    custom "Table_0" style to the built-in "No Style, Table Grid" style.
# 2) The presentation's design/theme colour scheme is switched from the
#    "Red Violet" (Integral) palette back to the default "Office" palette.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 5 -------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{098B196E-3227-442F-8B45-2A278DE237E1}")

# --- 2) Theme colour scheme: Red Violet -> Office -------------------------
$colors = $p.SlideMaster.Theme.ThemeColorScheme
$colors.Item(1).RGB  = 0          # dk1      000000
$colors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72
